# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Reorders/updates the worker-period detail table (rows 16-27, cols C:G)
# on sheet "Hoja1" so it is grouped by Periodo Mora (2210, 2302, 2303, 2304)
# and gives every worker a single, consistent Salario Basico.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Tipo Doc Trabajador | N Doc Trabajador | Nombre Trabajador | Periodo Mora | Valor Mora | Salario Basico
$data = @(
  @("CC", "1051448153", "AMIR PAJARO PAJARO",               "2210", 40000, 877803),
  @("CC", "1001898501", "DIDIER ANDRES MIRANDA SANCHEZ",    "2210", 40000, 1117172),
  @("CC", "1007856972", "HYLEANA MARGARITA BARRIOS PUERTA", "2210", 40000, 908526),
  @("CC", "1099962566", "MILTON JESUS CONDE LOZANO",        "2210", 40000, 877803),
  @("CC", "1001898501", "DIDIER ANDRES MIRANDA SANCHEZ",    "2302", 40000, 1117172),
  @("CC", "1007856972", "HYLEANA MARGARITA BARRIOS PUERTA", "2302", 46400, 908526),
  @("CC", "1099962566", "MILTON JESUS CONDE LOZANO",        "2302", 46400, 877803),
  @("CC", "1051448153", "AMIR PAJARO PAJARO",               "2303", 46400, 877803),
  @("CC", "1001898501", "DIDIER ANDRES MIRANDA SANCHEZ",    "2303", 46400, 1117172),
  @("CC", "1007856972", "HYLEANA MARGARITA BARRIOS PUERTA", "2303", 46400, 908526),
  @("CC", "1099962566", "MILTON JESUS CONDE LOZANO",        "2303", 46400, 877803),
  @("CC", "1007856972", "HYLEANA MARGARITA BARRIOS PUERTA", "2304", 32707, 908526)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $startRow + $i
  $rec = $data[$i]
  $ws.Cells.Item($row, 2).Value = $rec[0]
  $ws.Cells.Item($row, 3).Value = $rec[1]
  $ws.Cells.Item($row, 4).Value = $rec[2]
  $ws.Cells.Item($row, 5).Value = $rec[3]
  $ws.Cells.Item($row, 6).Value = $rec[4]
  $ws.Cells.Item($row, 7).Value = $rec[5]
}

# Columns were autosized ("best fit") in Excel after the data refresh.
$ws.Range("B:J").Columns.AutoFit() | Out-Null
